$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 100.833336
$ws.Range("I8").Value = 81
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 243
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = -104
$ws.Range("N8").Value = -878
$ws.Range("H17").Value = 409.66666
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 405.28088
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 1215.84264
$ws.Range("M17").Value = -2232
$ws.Range("N17").Value = -1551.84264
$ws.Range("H40").Value = 1998.625
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1998.625
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1998.625
$ws.Range("N40").Value = -2348.625
$ws.Range("M40").ClearContents()
$ws.Range("H64").Value = 3844.95
$ws.Range("I64").Value = 3824.9167
$ws.Range("J64").Value = 3875
$ws.Range("K64").Value = 3824.9167
$ws.Range("L64").Value = 3875
$ws.Range("M64").Value = -3576.9167
$ws.Range("N64").Value = -4371
$ws.Range("H67").Value = 3844.95
$ws.Range("I67").Value = 3824.9167
$ws.Range("J67").Value = 3875
$ws.Range("K67").Value = 3824.9167
$ws.Range("L67").Value = 3875
$ws.Range("M67").Value = -2966.9167
$ws.Range("N67").Value = -5591
$ws.Range("H127").Value = 2217
$ws.Range("I127").Value = 860
$ws.Range("J127").Value = 3260.8462
$ws.Range("K127").Value = 2580
$ws.Range("L127").Value = 9782.5386
$ws.Range("M127").Value = 2380
$ws.Range("N127").Value = -19702.5386
$ws.Range("H141").Value = 8274.666999999999
$ws.Range("I141").Value = 9267.916999999999
$ws.Range("J141").Value = 4301.6665
$ws.Range("K141").Value = 27803.751
$ws.Range("L141").Value = 12904.9995
$ws.Range("M141").Value = -22623.751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1306.4286
$ws.Range("I45").Value = 745
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 745
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -368
$ws.Range("N45").Value = -2154
$ws.Range("H61").Value = 5110.385
$ws.Range("I61").Value = 5385.6523
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 5385.6523
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -5173.6523
$ws.Range("N61").Value = -3424
$ws.Range("H132").Value = 2977373.5
$ws.Range("I132").Value = 3472963.5
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 10418890.5
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -10416360.5
$ws.Range("N132").Value = -16559
$ws.Range("H136").Value = 5110.385
$ws.Range("I136").Value = 5385.6523
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 16156.9569
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -13606.9569
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H107").Value = 33335034
$ws.Range("I107").Value = 52632860
$ws.Range("J107").Value = 2428.4546
$ws.Range("K107").Value = 52632860
$ws.Range("L107").Value = 2428.4546
$ws.Range("M107").Value = -52630940
$ws.Range("N107").Value = -6268.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8384.315000000001
$ws.Range("I31").Value = 939.05884
$ws.Range("J31").Value = 19390.348
$ws.Range("K31").Value = 939.05884
$ws.Range("L31").Value = 19390.348
$ws.Range("M31").Value = -644.05884
$ws.Range("N31").Value = -19980.348
$ws.Range("H34").Value = 8384.315000000001
$ws.Range("I34").Value = 939.05884
$ws.Range("J34").Value = 19390.348
$ws.Range("K34").Value = 939.05884
$ws.Range("L34").Value = 19390.348
$ws.Range("M34").Value = -737.05884
$ws.Range("N34").Value = -19794.348
$ws.Range("H58").Value = 4002519.2
$ws.Range("I58").Value = 7572217.5
$ws.Range("J58").Value = 12856.235
$ws.Range("K58").Value = 7572217.5
$ws.Range("L58").Value = 12856.235
$ws.Range("M58").Value = -7572014.5
$ws.Range("N58").Value = -13262.235
$ws.Range("H132").Value = 14498483
$ws.Range("I132").Value = 19608860
$ws.Range("J132").Value = 19085.666
$ws.Range("K132").Value = 58826580
$ws.Range("L132").Value = 57256.99800000001
$ws.Range("M132").Value = -58824050
$ws.Range("N132").Value = -62316.99800000001
$ws.Range("H134").Value = 8014143.5
$ws.Range("I134").Value = 10001189
$ws.Range("J134").Value = 4465848.5
$ws.Range("K134").Value = 30003567
$ws.Range("L134").Value = 13397545.5
$ws.Range("M134").Value = -30001032
$ws.Range("N134").Value = -13402615.5
$ws.Range("H136").Value = 4002519.2
$ws.Range("I136").Value = 7572217.5
$ws.Range("J136").Value = 12856.235
$ws.Range("K136").Value = 22716652.5
$ws.Range("L136").Value = 38568.705
$ws.Range("M136").Value = -22714102.5
$ws.Range("N136").Value = -43668.705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 420
$ws.Range("I10").Value = 56
$ws.Range("J10").Value = 875
$ws.Range("K10").Value = 168
$ws.Range("L10").Value = 2625
$ws.Range("M10").Value = -29
$ws.Range("N10").Value = -2903
$ws.Range("H104").Value = 58796.527
$ws.Range("I104").Value = 2022.5
$ws.Range("J104").Value = 73936.266
$ws.Range("K104").Value = 6067.5
$ws.Range("L104").Value = 221808.798
$ws.Range("M104").Value = -3446.5
$ws.Range("N104").Value = -227050.798
$ws.Range("H117").Value = 18471.666
$ws.Range("I117").Value = 33776.668
$ws.Range("J117").Value = 3166.6667
$ws.Range("K117").Value = 101330.004
$ws.Range("L117").Value = 9500.000100000001
$ws.Range("M117").Value = -97888.00399999999
$ws.Range("N117").Value = -16384.0001
$ws.Range("H129").Value = 1545.9615
$ws.Range("I129").Value = 506
$ws.Range("J129").Value = 1793.5714
$ws.Range("K129").Value = 1518
$ws.Range("L129").Value = 5380.7142
$ws.Range("M129").Value = 3482
$ws.Range("N129").Value = -15380.7142
$ws.Range("H134").Value = 3697.75
$ws.Range("I134").Value = 3000.238
$ws.Range("J134").Value = 5029.364
$ws.Range("K134").Value = 9000.714
$ws.Range("L134").Value = 15088.092
$ws.Range("M134").Value = -3930.714
$ws.Range("N134").Value = -25228.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1502.0834
$ws.Range("I7").Value = 1339.3572
$ws.Range("J7").Value = 1729.9
$ws.Range("K7").Value = 1339.3572
$ws.Range("L7").Value = 1729.9
$ws.Range("M7").Value = -1227.3572
$ws.Range("H126").Value = 1502.0834
$ws.Range("I126").Value = 1339.3572
$ws.Range("J126").Value = 1729.9
$ws.Range("K126").Value = 4018.0716
$ws.Range("L126").Value = 5189.700000000001
$ws.Range("M126").Value = -1548.0716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 155827420
$ws.Range("I132").Value = 171436430
$ws.Range("J132").Value = 128511660
$ws.Range("K132").Value = 514309290
$ws.Range("L132").Value = 385534980
$ws.Range("M132").Value = -514306760
$ws.Range("N132").Value = -385540040
